$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-04-29 Monday" "2024-04-30 Tuesday"

Replace-Text "573×7=4011" "682×3=2046"
Replace-Text "648×5=3240" "282×3=846"
Replace-Text "301×3=903" "911×5=4555"
Replace-Text "946×9=8514" "767×8=6136"
Replace-Text "797×6=4782" "762×7=5334"

Replace-Text "255×9=2295" "828×4=3312"
Replace-Text "868×8=6944" "631×2=1262"
Replace-Text "624×3=1872" "354×6=2124"
Replace-Text "538×5=2690" "657×9=5913"
Replace-Text "758×3=2274" "508×7=3556"

Replace-Text "518×2=1036" "177×9=1593"
Replace-Text "272×3=816" "902×4=3608"
Replace-Text "706×2=1412" "241×4=964"
Replace-Text "865×8=6920" "442×9=3978"
Replace-Text "305×2=610" "238×3=714"

Replace-Text "291×7=2037" "743×8=5944"
Replace-Text "917×3=2751" "508×6=3048"
Replace-Text "660×5=3300" "999×6=5994"
Replace-Text "731×2=1462" "256×6=1536"
Replace-Text "834×5=4170" "384×3=1152"

Replace-Text "501×4=2004" "556×2=1112"
Replace-Text "194×3=582" "569×7=3983"
Replace-Text "537×9=4833" "693×7=4851"
Replace-Text "287×2=574" "877×7=6139"
Replace-Text "200×8=1600" "465×8=3720"
